$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update installed power for "Motores Elétricos" (Electric Motors) row
$ws.Range("B6").Value = 22.69

# Force recalculation so dependent formulas (D6, D7) update
$excel.Calculate()

# Update the active selection to match the final state
$ws.Range("B7").Select()
